$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 23.699655
$ws.Range("H2").Value = 71.09896499999999
$ws.Range("I2").Value = 0.4841969272415696
$ws.Range("J2").Value = 0.4841969272415697
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 1446.665010745815
$ws.Range("R2").Value = 13019.98509671234
$ws.Range("S2").Value = 0.09895113580797001
$ws.Range("T2").Value = 0.09895113580797002
$ws.Range("G3").Value = 23.699655
$ws.Range("H3").Value = 71.09896499999999
$ws.Range("I3").Value = 0.4841969272415696
$ws.Range("J3").Value = 0.4841969272415697
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 2519.616165709229
$ws.Range("R3").Value = 22676.54549138307
$ws.Range("S3").Value = 0.1723404378657893
$ws.Range("T3").Value = 0.1723404378657894
$ws.Range("G4").Value = 23.699655
$ws.Range("H4").Value = 71.09896499999999
$ws.Range("I4").Value = 0.4841969272415696
$ws.Range("J4").Value = 0.4841969272415697
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 3112.674989448779
$ws.Range("R4").Value = 28014.07490503901
$ws.Range("S4").Value = 0.2129053535678102
$ws.Range("T4").Value = 0.2129053535678103
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3952142927098025
$ws.Range("J5").Value = 0.3952142927098025
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 1180.806107686592
$ws.Range("R5").Value = 10627.25496917933
$ws.Range("S5").Value = 0.08076652483931963
$ws.Range("T5").Value = 0.08076652483931963
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3952142927098025
$ws.Range("J6").Value = 0.3952142927098025
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.140668807306256
$ws.Range("T6").Value = 0.140668807306256
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3952142927098025
$ws.Range("J7").Value = 0.3952142927098025
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 2540.647358913841
$ws.Range("R7").Value = 22865.82623022457
$ws.Range("S7").Value = 0.1737789605642268
$ws.Range("T7").Value = 0.1737789605642269
$ws.Range("G8").Value = 5.902376333333333
$ws.Range("H8").Value = 17.707129
$ws.Range("I8").Value = 0.1205887800486278
$ws.Range("J8").Value = 0.1205887800486278
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 360.2905325705168
$ws.Range("R8").Value = 3242.614793134651
$ws.Range("S8").Value = 0.02464368541016377
$ws.Range("T8").Value = 0.02464368541016377
$ws.Range("G9").Value = 5.902376333333333
$ws.Range("H9").Value = 17.707129
$ws.Range("I9").Value = 0.1205887800486278
$ws.Range("J9").Value = 0.1205887800486278
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 627.5079880093713
$ws.Range("R9").Value = 5647.571892084342
$ws.Range("S9").Value = 0.04292122065639095
$ws.Range("T9").Value = 0.04292122065639095
$ws.Range("G10").Value = 5.902376333333333
$ws.Range("H10").Value = 17.707129
$ws.Range("I10").Value = 0.1205887800486278
$ws.Range("J10").Value = 0.1205887800486278
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 775.2087189067123
$ws.Range("R10").Value = 6976.878470160411
$ws.Range("S10").Value = 0.05302387398207311
$ws.Range("T10").Value = 0.05302387398207312